# Append a new data row (row 6) to the order inventory sheet, mirroring the
# text-formatted values used by the existing rows (order #, part #, item,
# price, seller, buyer).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new row to be stored as text (like all the other data rows),
# so numeric-looking values such as "30" don't get reinterpreted as numbers.
$ws.Range("A6:F6").NumberFormat = "@"

$ws.Range("A6").Value = "30"
$ws.Range("B6").Value = "23402043"
$ws.Range("C6").Value = "book"
$ws.Range("D6").Value = "10"
$ws.Range("E6").Value = "Eric"
$ws.Range("F6").Value = "Fred"
